$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of (accuracy header column, min column, max column)
$pairs = @(
    @{ Header = "B"; Min = "C"; Max = "D" },
    @{ Header = "E"; Min = "F"; Max = "G" },
    @{ Header = "H"; Min = "I"; Max = "J" },
    @{ Header = "K"; Min = "L"; Max = "M" },
    @{ Header = "N"; Min = "O"; Max = "P" },
    @{ Header = "Q"; Min = "R"; Max = "S" },
    @{ Header = "T"; Min = "U"; Max = "V" }
)

foreach ($pair in $pairs) {
    $ws.Range($pair.Min + "1").Value = "min"
    $ws.Range($pair.Max + "1").Value = "max"
}

$ws.Range("N13").Select()
